$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 66; $r++) {
    if ($r -eq 46) {
        $ws.Cells.Item($r, 6).Value = 3
    } else {
        $ws.Cells.Item($r, 6).Value = 2
    }
}
